$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.45845575611036
$ws.Range("D2").Value = 3.217835152341511
$ws.Range("E2").Value = 25.36626754610605
$ws.Range("F2").Value = 17.25248603758778
$ws.Range("G2").Value = 17.36485809358217
$ws.Range("H2").Value = 10.95481938028163
$ws.Range("I2").Value = 24.88978878512567
$ws.Range("L2").Value = 9.478113825664025
$ws.Range("M2").Value = 15.28720996369622
$ws.Range("O2").Value = 15.16812747536079
$ws.Range("B3").Value = 15.94367005869347
$ws.Range("D3").Value = 3.159650268412608
$ws.Range("E3").Value = 24.94554764005377
$ws.Range("F3").Value = 17.16108644283293
$ws.Range("G3").Value = 17.11070209341355
$ws.Range("H3").Value = 10.987134010166
$ws.Range("I3").Value = 24.89207753931743
$ws.Range("L3").Value = 9.297033254702527
$ws.Range("M3").Value = 15.02945467375484
$ws.Range("O3").Value = 15.17577376810196
$ws.Range("B4").Value = 15.6183845316253
$ws.Range("D4").Value = 3.123655407579521
$ws.Range("E4").Value = 24.68717998602125
$ws.Range("F4").Value = 17.11181683039942
$ws.Range("G4").Value = 16.96285270803779
$ws.Range("H4").Value = 11.00963196508436
$ws.Range("I4").Value = 24.90163803809496
$ws.Range("L4").Value = 9.184129021331568
$ws.Range("M4").Value = 14.86896751404717
$ws.Range("O4").Value = 15.18595763283854
$ws.Range("B5").Value = 15.48367509447023
$ws.Range("D5").Value = 3.108937139290828
$ws.Range("E5").Value = 24.58202147833993
$ws.Range("F5").Value = 17.09348033306944
$ws.Range("G5").Value = 16.90476455959222
$ws.Range("H5").Value = 11.01946614621141
$ws.Range("I5").Value = 24.90758857043684
$ws.Range("L5").Value = 9.137735977347802
$ws.Range("M5").Value = 14.80307157477163
$ws.Range("O5").Value = 15.19148382554263
$ws.Range("B6").Value = 15.46118166732514
$ws.Range("D6").Value = 3.106490680542723
$ws.Range("E6").Value = 24.56457186158082
$ws.Range("F6").Value = 17.09054122041085
$ws.Range("G6").Value = 16.89525238193869
$ws.Range("H6").Value = 11.02113927498065
$ws.Range("I6").Value = 24.90870082521666
$ws.Range("L6").Value = 9.13001073780341
$ws.Range("M6").Value = 14.79210152139856
$ws.Range("O6").Value = 15.19248442248577
$ws.Range("B7").Value = 15.61657628615022
$ws.Range("D7").Value = 3.123457091025138
$ws.Range("E7").Value = 24.68576107584607
$ws.Range("F7").Value = 17.11156246544827
$ws.Range("G7").Value = 16.96206043281179
$ws.Range("H7").Value = 11.00976189802415
$ws.Range("I7").Value = 24.90170996627763
$ws.Range("L7").Value = 9.183504836094604
$ws.Range("M7").Value = 14.86808074235599
$ws.Range("O7").Value = 15.18602659513194
$ws.Range("B8").Value = 16.28296080093604
$ws.Range("D8").Value = 3.197839275877797
$ws.Range("E8").Value = 25.22129536618107
$ws.Range("F8").Value = 17.21956202061706
$ws.Range("G8").Value = 17.27558380340668
$ws.Range("H8").Value = 10.96540893736773
$ws.Range("I8").Value = 24.88888810073969
$ws.Range("L8").Value = 9.416062997734812
$ws.Range("M8").Value = 15.19883496593612
$ws.Range("O8").Value = 15.16962257032159
$ws.Range("B9").Value = 17.51031211221089
$ws.Range("D9").Value = 3.342430178855475
$ws.Range("E9").Value = 26.26513806046673
$ws.Range("F9").Value = 17.48471103094569
$ws.Range("G9").Value = 17.95084327298596
$ws.Range("H9").Value = 10.89959508279119
$ws.Range("I9").Value = 24.9282159859705
$ws.Range("L9").Value = 9.856410184077426
$ws.Range("M9").Value = 15.8271163073317
$ws.Range("O9").Value = 15.18113561121295
$ws.Range("B10").Value = 18.35604300654596
$ws.Range("D10").Value = 3.464578552378324
$ws.Range("E10").Value = 27.02041974202161
$ws.Range("F10").Value = 17.7104941113035
$ws.Range("G10").Value = 18.47726366782242
$ws.Range("H10").Value = 10.86424935747377
$ws.Range("I10").Value = 24.99599871058021
$ws.Range("L10").Value = 10.16776475049617
$ws.Range("M10").Value = 16.27289668302444
$ws.Range("O10").Value = 15.21631992442333
$ws.Range("B11").Value = 18.72728648971335
$ws.Range("D11").Value = 3.518298364037961
$ws.Range("E11").Value = 27.35993920836514
$ws.Range("F11").Value = 17.81954292777671
$ws.Range("G11").Value = 18.72189075053835
$ws.Range("H11").Value = 10.85101438382307
$ws.Range("I11").Value = 25.03515809588107
$ws.Range("L11").Value = 10.30626124219064
$ws.Range("M11").Value = 16.47158323997095
$ws.Range("O11").Value = 15.23812363490352
$ws.Range("B12").Value = 18.86583212980472
$ws.Range("D12").Value = 3.538363517851526
$ws.Range("E12").Value = 27.48780031477444
$ws.Range("F12").Value = 17.86171172179489
$ws.Range("G12").Value = 18.8151461826879
$ws.Range("H12").Value = 10.84641302807382
$ws.Range("I12").Value = 25.0511711269417
$ws.Range("L12").Value = 10.35821581572247
$ws.Range("M12").Value = 16.54617809160563
$ws.Range("O12").Value = 15.24721128474586
$ws.Range("B13").Value = 18.83608583861316
$ws.Range("D13").Value = 3.534054642883467
$ws.Range("E13").Value = 27.46029650541242
$ws.Range("F13").Value = 17.85259164420126
$ws.Range("G13").Value = 18.79503626944468
$ws.Range("H13").Value = 10.8473857359262
$ws.Range("I13").Value = 25.04766998670164
$ws.Range("L13").Value = 10.3470489096536
$ws.Range("M13").Value = 16.53014219628563
$ws.Range("O13").Value = 15.2452171924365
$ws.Range("B14").Value = 18.7387260490379
$ws.Range("D14").Value = 3.519954759837437
$ws.Range("E14").Value = 27.37047328380441
$ws.Range("F14").Value = 17.82299487666794
$ws.Range("G14").Value = 18.72955118298199
$ws.Range("H14").Value = 10.85062759372282
$ws.Range("I14").Value = 25.03645184313472
$ws.Range("L14").Value = 10.31054563457793
$ws.Range("M14").Value = 16.47773339582996
$ws.Range("O14").Value = 15.23885464591047
$ws.Range("B15").Value = 18.67882249645134
$ws.Range("D15").Value = 3.511281716370767
$ws.Range("E15").Value = 27.31535826574802
$ws.Range("F15").Value = 17.8049787310162
$ws.Range("G15").Value = 18.68951691709623
$ws.Range("H15").Value = 10.85266681815844
$ws.Range("I15").Value = 25.02973421900984
$ws.Range("L15").Value = 10.28812124789585
$ws.Range("M15").Value = 16.44554615449816
$ws.Range("O15").Value = 15.23506552770298
$ws.Range("B16").Value = 18.3315014721472
$ws.Range("D16").Value = 3.461029646783156
$ws.Range("E16").Value = 26.99813891670083
$ws.Range("F16").Value = 17.70349194513788
$ws.Range("G16").Value = 18.46137009084963
$ws.Range("H16").Value = 10.86517165231743
$ws.Range("I16").Value = 24.99360600278537
$ws.Range("L16").Value = 10.15864705685344
$ws.Range("M16").Value = 16.25982486682545
$ws.Range("O16").Value = 15.2150114730521
$ws.Range("B17").Value = 18.11490494643174
$ws.Range("D17").Value = 3.42971999958272
$ws.Range("E17").Value = 26.8024049854324
$ws.Range("F17").Value = 17.64283092928868
$ws.Range("G17").Value = 18.32264348143622
$ws.Range("H17").Value = 10.87357247581541
$ws.Range("I17").Value = 24.97356600777152
$ws.Range("L17").Value = 10.07838565132662
$ws.Range("M17").Value = 16.14480133863291
$ws.Range("O17").Value = 15.20419250202293
$ws.Range("B18").Value = 17.98906026697715
$ws.Range("D18").Value = 3.411538400230288
$ws.Range("E18").Value = 26.68944867900328
$ws.Range("F18").Value = 17.60853934708107
$ws.Range("G18").Value = 18.24334408416242
$ws.Range("H18").Value = 10.87867201651299
$ws.Range("I18").Value = 24.96282396006842
$ws.Range("L18").Value = 10.03192803885728
$ws.Range("M18").Value = 16.07825985422601
$ws.Range("O18").Value = 15.19851571528176
$ws.Range("B19").Value = 17.9462374621796
$ws.Range("D19").Value = 3.405353069595962
$ws.Range("E19").Value = 26.65114301368097
$ws.Range("F19").Value = 17.59703278588529
$ws.Range("G19").Value = 18.21658296128722
$ws.Range("H19").Value = 10.88044454787185
$ws.Range("I19").Value = 24.95932197986372
$ws.Range("L19").Value = 10.01614914623512
$ws.Range("M19").Value = 16.05566595523015
$ws.Range("O19").Value = 15.19668748990824
$ws.Range("B20").Value = 18.13809358479436
$ws.Range("D20").Value = 3.433070971592221
$ws.Range("E20").Value = 26.82328095515287
$ws.Range("F20").Value = 17.64922666277972
$ws.Range("G20").Value = 18.33736107672328
$ws.Range("H20").Value = 10.87265048603018
$ws.Range("I20").Value = 24.97561820237939
$ws.Range("L20").Value = 10.08696027210168
$ws.Range("M20").Value = 16.15708582610981
$ws.Range("O20").Value = 15.20528770655608
$ws.Range("B21").Value = 18.76737896123628
$ws.Range("D21").Value = 3.524103854896779
$ws.Range("E21").Value = 27.39687669690962
$ws.Range("F21").Value = 17.83166474140726
$ws.Range("G21").Value = 18.74876985679895
$ws.Range("H21").Value = 10.84966423120039
$ws.Range("I21").Value = 25.03971485449338
$ws.Range("L21").Value = 10.32128116014552
$ws.Range("M21").Value = 16.49314500693752
$ws.Range("O21").Value = 15.24070095353487
$ws.Range("B22").Value = 19.16674732828545
$ws.Range("D22").Value = 3.581977557268734
$ws.Range("E22").Value = 27.76757773355791
$ws.Range("F22").Value = 17.95597638706202
$ws.Range("G22").Value = 19.02121110159419
$ws.Range("H22").Value = 10.83703432876688
$ws.Range("I22").Value = 25.08850123037642
$ws.Range("L22").Value = 10.47154470877203
$ws.Range("M22").Value = 16.70900585038768
$ws.Range("O22").Value = 15.2686874135473
$ws.Range("B23").Value = 18.95471560170349
$ws.Range("D23").Value = 3.551241315186404
$ws.Range("E23").Value = 27.57014896598915
$ws.Range("F23").Value = 17.8891771254911
$ws.Range("G23").Value = 18.87551765865456
$ws.Range("H23").Value = 10.84355574179494
$ws.Range("I23").Value = 25.06183672736702
$ws.Range("L23").Value = 10.39162201476739
$ws.Range("M23").Value = 16.59415908670958
$ws.Range("O23").Value = 15.25330871910899
$ws.Range("B24").Value = 18.12761410678644
$ws.Range("D24").Value = 3.431556560834392
$ws.Range("E24").Value = 26.8138442474695
$ws.Range("F24").Value = 17.64633333210102
$ws.Range("G24").Value = 18.33070582108011
$ws.Range("H24").Value = 10.87306647731511
$ws.Range("I24").Value = 24.97468797704656
$ws.Range("L24").Value = 10.08308466187233
$ws.Range("M24").Value = 16.15153329136439
$ws.Range("O24").Value = 15.20479087240159
$ws.Range("B25").Value = 17.18761916884394
$ws.Range("D25").Value = 3.302549633918447
$ws.Range("E25").Value = 25.98424847407998
$ws.Range("F25").Value = 17.40741671513138
$ws.Range("G25").Value = 17.76236082560722
$ws.Range("H25").Value = 10.9151223891444
$ws.Range("I25").Value = 24.91071126444923
$ws.Range("L25").Value = 9.739242213952446
$ws.Range("M25").Value = 15.65968887962068
$ws.Range("O25").Value = 15.17332973051478
